$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Frasca Luca"
$ws.Range("B3").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C3").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("D3").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E3").Value = "daniel pedrotti | iMontagna"
$ws.Range("F3").Value = "Davide Raffaelli | MediaserT"
